# Trade #7 closed at 2026-02-16 22:52:38 - base_strategy DOWN +0.000%
# Append a new trade row (row 8) to both the "All Trades" and
# "base_strategy" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(8, 1).Value = 7

    # "2026-02-16" looks like a date, so force the cell to text first,
    # write the literal string, then drop the format override again so
    # the stored value matches the plain text used by the rows above it.
    $ws.Cells.Item(8, 2).NumberFormat = "@"
    $ws.Cells.Item(8, 2).Value = "2026-02-16"
    $ws.Cells.Item(8, 2).Style = "Normal"

    $ws.Cells.Item(8, 3).Value = "22:52:38"
    $ws.Cells.Item(8, 4).Value = "base_strategy"
    $ws.Cells.Item(8, 5).Value = "DOWN"
    $ws.Cells.Item(8, 6).Value = 49.999998
    $ws.Cells.Item(8, 7).Value = ""
    $ws.Cells.Item(8, 8).Value = "OPEN"
    $ws.Cells.Item(8, 9).Value = 0
    $ws.Cells.Item(8, 10).Value = 0
    $ws.Cells.Item(8, 11).Value = 100
    $ws.Cells.Item(8, 12).Value = 0
    $ws.Cells.Item(8, 13).Value = 0
    $ws.Cells.Item(8, 14).Value = 0.6
    $ws.Cells.Item(8, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(8, 16).Value = ""
    $ws.Cells.Item(8, 17).Value = 0
}
